# Generate Report for Handback
# Update the timestamp cells to reflect the new handoff/handback generation times.

$wb = $excel.ActiveWorkbook

# Sheet "Overview": G2 - Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-24 13:10:45"

# Sheet "zh-cn": H2 - Correspond Handoff Datetime, K2 - Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-24 13:10:40"
$wsZhCn.Range("K2").Value = "2016-08-24 13:10:58"

# Sheet "de-de": H2 - Correspond Handoff Datetime, K2 - Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-24 13:10:45"
$wsDeDe.Range("K2").Value = "2016-08-24 13:11:23"
